$d = $word.ActiveDocument

function Replace-Paragraph($findText, $newXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    $para = $rng.Paragraphs(1).Range
    $para.InsertXML($newXml)
}

Replace-Paragraph 'Отчет по тестированию' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:b w:val="on"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="64"/></w:rPr><w:t>Test Report</w:t></w:r></w:p>'
Replace-Paragraph 'Результаты тестирования за 03 февраля 2023' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="32"/><w:position w:val="20"/><w:u w:val="single"/></w:rPr><w:t>Testing result for date 18 апреля 2023</w:t></w:r></w:p>'
Replace-Paragraph 'Всего: 4; ' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t xml:space="preserve">Total: 4; </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t xml:space="preserve">Successful: 0; </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t xml:space="preserve">Failed: 4; </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t xml:space="preserve">Broken: 0; </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t xml:space="preserve">Skipped: 0; </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:position w:val="20"/><w:b w:val="on"/></w:rPr><w:t>Unknown: 0.</w:t></w:r></w:p>'
Replace-Paragraph 'Тест кейс' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:b w:val="on"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>Test case</w:t></w:r></w:p>'
Replace-Paragraph 'Статус' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:b w:val="on"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>Status</w:t></w:r></w:p>'
Replace-Paragraph 'Тест кейсы' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:b w:val="on"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="64"/></w:rPr><w:t>Test cases</w:t></w:r></w:p>'
